$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the authoritative diff.
# Numeric-looking text (e.g. "2.50", "1.00") must stay TEXT, matching the
# original inlineStr cells -- Excel auto-coerces such strings to numbers on
# a plain Range.Value assignment, so those cells are forced to a text
# number-format first, then the style is reset back to "Normal" so no stray
# formatting is introduced (the source cells carry no explicit style).

$ws.Range("D2").Value = "69.441.44"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.769.07"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "3.766.19"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.484"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "4.396.42"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "3.769.44"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "69.521.03"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.730"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000135"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.138"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.341"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "457.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").Value = "2.956.65"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0360"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
